# Fill in the newly-tracked weekly figures (columns B..G) for rows 105-125
# on the active sheet, then move the view/selection the way the author left
# it (top-left cell A79, active cell K108).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(B, C, D, E, F, G)  -- $null means "leave that cell blank"
$rows = [ordered]@{
    105 = @(32, 266, 44, 5226, 145, $null)
    106 = @(36, 321, 58, 5211, 77, 119)
    107 = @(40, 321, 54, 5230, 138, 79)
    108 = @(37, 451, 90, $null, $null, $null)
    109 = @(19, 296, 55, 5297, 165, 41)
    110 = @(19, 444, 87, 5294, 140, 88)
    111 = @(12, 324, 43, 5350, 138, 46)
    112 = @(5, 216, 34, 5326, 108, 105)
    113 = @(11, 272, 44, 5368, 74, 108)
    114 = @(24, 346, 61, 5430, 128, 16)
    115 = @(23, 311, 48, 5429, 144, 24)
    116 = @(10, 229, 39, 5389, 62, 156)
    117 = @(0, 197, 32, 5433, 76, 98)
    118 = @(2, 238, 39, $null, $null, $null)
    119 = @(3, 153, 22, 5488, 82, 42)
    120 = @(7, 197, 31, 5514, 80, 25)
    121 = @(9, 170, 26, 5528, 49, 51)
    122 = @(0, 166, 28, $null, $null, $null)
    123 = @(10, 189, 28, $null, $null, $null)
    124 = @(0, 149, 23, 5523, 115, $null)
    125 = @(8, 202, 34, 5529, 84, 33)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $v = $values[$i]
        if ($null -ne $v) {
            $ws.Range($columns[$i] + $r).Value = $v
        }
    }
}

# Move the selection/scroll position to match where the author ended up.
$ws.Range("K108").Select()
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 1
